$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their text/string representation (avoid numeric auto-conversion)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "36.790.56"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.126.55"
$ws.Range("E3").Value = "  +10.43%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "256.20"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "0.668"
$ws.Range("E6").Value = "  -4.58%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "47.48"
$ws.Range("E8").Value = "  +7.06%  "
$ws.Range("D9").Value = "61.14"
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "2.435.49"
$ws.Range("E13").Value = "  +10.60%  "
$ws.Range("D14").Value = "14.50"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("E15").Value = "  +5.57%  "
$ws.Range("D16").Value = "2.125.12"
$ws.Range("E16").Value = "  +10.22%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "36.936.74"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "73.94"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "0.0₃0842"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "242.22"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("D23").Value = "5.24"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E25").Value = "  -7.79%  "
$ws.Range("D26").Value = "171.90"
$ws.Range("D27").Value = "21.58"
$ws.Range("E27").Value = "  +14.52%  "
$ws.Range("E28").Value = "  +4.35%  "
$ws.Range("E29").Value = "  -9.47%  "
$ws.Range("D30").Value = "28.02"
$ws.Range("E30").Value = "  +60.51%  "
$ws.Range("E31").Value = "  -4.69%  "
$ws.Range("D32").Value = "4.53"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "0.0955"
$ws.Range("E33").Value = "  +11.08%  "
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("E35").Value = "  +16.45%  "
$ws.Range("D36").Value = "0.948"
$ws.Range("E36").Value = "  +9.45%  "
$ws.Range("E37").Value = "  -4.63%  "
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("E40").Value = "  -8.01%  "
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").Value = "99.49"
$ws.Range("E43").Value = "  -8.90%  "
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  +13.52%  "
$ws.Range("D45").Value = "16.36"
$ws.Range("E45").Value = "  -5.42%  "
$ws.Range("D46").Value = "1.358.62"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("E47").Value = "  +3.86%  "
$ws.Range("D48").Value = "7.13"
$ws.Range("E48").Value = "  +10.41%  "
$ws.Range("D49").Value = "2.309.21"
$ws.Range("E49").Value = "  +9.88%  "
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("D51").Value = "2.84"
$ws.Range("E51").Value = "  +0.83%  "
